$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest crypto snapshot.
# D-column numeric-looking text is entered with a leading apostrophe so Excel
# keeps it as literal text (preserving formats like '1.110' or '78.10'), then
# the style is reset to Normal so no extra quote-prefix styling is persisted.
$ws.Range("D2").Value = "'27.997.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "'1.858.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'312.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "'0.5139"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").Value = "'0.08218"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.15%  "

$ws.Range("D10").Value = "'1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").Value = "'41.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").Value = "'6.193"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.61%  "

$ws.Range("D13").Value = "'20.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").Value = "'1.862.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").Value = "'7.259"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").Value = "'0.00001097"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").Value = "'90.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "'17.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "'6.014"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").Value = "'28.024.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("E24").Value = "  -3.50%  "

$ws.Range("D25").Value = "'2.260"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").Value = "'2.073.23"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'2.509"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "

$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "'20.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").Value = "'124.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("D31").Value = "'0.1069"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").Value = "'1.032"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("D33").Value = "'5.912"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.11%  "

$ws.Range("D34").Value = "'3.594"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").Value = "'9.418"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.98%  "

$ws.Range("D36").Value = "'0.02411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("D38").Value = "'0.2181"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").Value = "'0.6556"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.04%  "

$ws.Range("D40").Value = "'1.196"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").Value = "'4.992"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("D42").Value = "'1.211"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("D43").Value = "'11.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.56%  "

$ws.Range("D44").Value = "'0.6148"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "

$ws.Range("D45").Value = "'12.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("D46").Value = "'1.281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").Value = "'3.659"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").Value = "'2.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "'120.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").Value = "'78.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.38%  "
